$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45/46/47 text is set in a particular order so that the shared
# string table ends up with the same slot assignment as the reference
# workbook: A45 is renamed first (mutating the formerly-unique string in
# place), then A46 re-uses that same text, then A47 gets a brand-new
# string, and finally A45 is renamed again to its real final text (which,
# since the slot is now shared with A46, becomes a new table entry). ---
$ws.Range("A45").Value = "Implementer Afskrivning"
$ws.Range("A46").Value = "Implementer Afskrivning"
$ws.Range("A47").Value = "OpdatererAfskrivningUi"
$ws.Range("A45").Value = "Opdatere AfskrivningUI"

# --- Row 45: update role (date/time untouched) ---
$ws.Range("B45").Value = "User-Interface Designer"

# --- Row 46: was a placeholder empty row, now filled in ---
$ws.Range("B46").Value = "Implementer"
$ws.Range("C46").Value = "2020-03-10"
$ws.Range("D46").Value = 0.375
$ws.Range("E46").Value = 0.45833333333333331
$ws.Range("G46").Formula = "=E46-D46"

# --- Row 47: was a placeholder empty row, now filled in ---
$ws.Range("B47").Value = "User-Interface Designer"
$ws.Range("C47").Value = "2020-03-10"
$ws.Range("D47").Value = 0.5
$ws.Range("E47").Value = 0.625
$ws.Range("G47").Formula = "=E47-D47"

# --- Rows 48-55: new rows. Column C on these rows uses a "date, centered
# (no vertical-centering)" style that doesn't otherwise exist on a blank
# cell, so seed the style by copying an existing date cell's format first,
# then overwrite with the real date value. ---
$ws.Range("C44").Copy()
$ws.Range("C48:C55").PasteSpecial(-4122)
$ws.Range("C48:C55").VerticalAlignment = -4107

# --- Row 48: new row ---
$ws.Range("A48").Value = "Tilføj scroll til grund UI"
$ws.Range("B48").Value = "User-Interface Designer"
$ws.Range("C48").Value = "2020-03-11"
$ws.Range("D48").Value = 0.375
$ws.Range("E48").Value = 0.3888888888888889
$ws.Range("G48").Formula = "=E48-D48"

# --- Row 49: new row ---
$ws.Range("A49").Value = "Reviewer AD09"
$ws.Range("B49").Value = "Reviewer"
$ws.Range("C49").Value = "2020-03-11"
$ws.Range("D49").Value = 0.39583333333333331
$ws.Range("E49").Value = 0.40972222222222227
$ws.Range("G49").Formula = "=E49-D49"

# --- Row 50: new row ---
$ws.Range("A50").Value = "Lav SSD09"
$ws.Range("B50").Value = "System Analyst "
$ws.Range("C50").Value = "2020-03-11"
$ws.Range("D50").Value = 0.41666666666666669
$ws.Range("E50").Value = 0.4375
$ws.Range("G50").Formula = "=E50-D50"

# --- Row 51: new row ---
$ws.Range("A51").Value = "Lav OC0903"
$ws.Range("B51").Value = "System Analyst "
$ws.Range("C51").Value = "2020-03-11"
$ws.Range("D51").Value = 0.4375
$ws.Range("E51").Value = 0.45833333333333331
$ws.Range("G51").Formula = "=E51-D51"

# --- Row 52: new row ---
$ws.Range("A52").Value = "Lav SD09"
$ws.Range("B52").Value = "Software Architect"
$ws.Range("C52").Value = "2020-03-11"
$ws.Range("D52").Value = 0.47916666666666669
$ws.Range("E52").Value = 0.5
$ws.Range("G52").Formula = "=E52-D52"

# --- Row 53: new row ---
$ws.Range("A53").Value = "Lav DCD09"
$ws.Range("B53").Value = "Software Architect"
$ws.Range("C53").Value = "2020-03-11"
$ws.Range("D53").Value = 0.52083333333333337
$ws.Range("E53").Value = 0.54166666666666663
$ws.Range("G53").Formula = "=E53-D53"

# --- Row 54: new row ---
$ws.Range("A54").Value = "Implementer UC09"
$ws.Range("B54").Value = "Implementer"
$ws.Range("C54").Value = "2020-03-11"
$ws.Range("D54").Value = 0.5625
$ws.Range("E54").Value = 0.60416666666666663
$ws.Range("G54").Formula = "=E54-D54"

# --- Row 55: new row ---
$ws.Range("A55").Value = "Lav SSD10"
$ws.Range("B55").Value = "System Analyst "
$ws.Range("C55").Value = "2020-03-11"
$ws.Range("D55").Value = 0.60416666666666663
$ws.Range("E55").Value = 0.625
$ws.Range("G55").Formula = "=E55-D55"

# --- Row 56: new row (note: no date value in column C for this one) ---
$ws.Range("A56").Value = "Lav OC10"
$ws.Range("B56").Value = "System Analyst "
$ws.Range("D56").Value = 0.64583333333333337
$ws.Range("E56").Value = 0.66666666666666663
$ws.Range("G56").Formula = "=E56-D56"

# Update selection / view to match the saved state (no frozen scroll position, select A9)
$ws.Range("A9").Select()
